$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Хід роботи" is followed by a single empty paragraph (bold, sz28,
#    spacing -67). Replace it with three paragraphs:
#      - an empty bold/sz28 paragraph (no -67 spacing)
#      - a new paragraph with the GitHub repository link
#      - the original empty paragraph (bold, sz28, spacing -67), now also
#        carrying the _GoBack bookmark
# ---------------------------------------------------------------------------
$hidRobotyPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "Хід роботи") {
        $hidRobotyPara = $cand
        break
    }
}
$target1 = $hidRobotyPara.Next()
$xml1 = @'
<w:p><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:line="322" w:lineRule="exact"/><w:ind w:left="440"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:line="322" w:lineRule="exact"/><w:ind w:left="440"/><w:rPr><w:b/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>GitHub</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:lang w:val="ru-RU"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">репозиторій: </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>https://github.com/AlexanderHorielko/SAI_Horielko_PI-59</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:line="322" w:lineRule="exact"/><w:ind w:left="440"/><w:rPr><w:b/><w:spacing w:val="-67"/><w:sz w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$target1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2. The "Рисунок 4. Модель кластеризації" caption paragraph previously held
#    the _GoBack bookmark between "Модель " and "кластеризації". Drop the
#    bookmark from here (it now lives on the new empty paragraph above) and
#    merge the trailing two runs into a single " кластеризації" run.
# ---------------------------------------------------------------------------
$figPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "Рисунок 4. Модель кластеризації") {
        $figPara = $cand
        break
    }
}
$xml2 = @'
<w:p><w:pPr><w:pStyle w:val="TableParagraph"/><w:spacing w:before="77"/><w:ind w:right="21"/><w:jc w:val="center"/><w:rPr><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>Рисунок</w:t></w:r><w:r><w:rPr><w:spacing w:val="-5"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="ru-RU"/></w:rPr><w:t>4</w:t></w:r><w:r><w:rPr><w:spacing w:val="-3"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:lang w:val="ru-RU"/></w:rPr><w:t>М</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>одель</w:t></w:r><w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> кластеризації</w:t></w:r></w:p>
'@
$figPara.Range.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3. The (empty) paragraph right after the caption switches its run-language
#    from ru-RU to en-US.
# ---------------------------------------------------------------------------
$figPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13) -eq "Рисунок 4. Модель кластеризації") {
        $figPara2 = $cand
        break
    }
}
$langPara = $figPara2.Next()
$langPara.Range.LanguageID = "en-US"

Write-Output "Done"
